# Apply the edits described by the commit:
#  - Fill in grade values (5) for row 5 (columns C:H) and add a new
#    graded cell in J5 (style copied from a same-styled neighbour, J15).
#  - Fill in grade values (5) for the previously-empty K9 cell and add a
#    new graded cell in J9 (both using the "border-right" style already
#    used on row 9's I cell family, mirrored from J21/K21).
#  - Re-point the frozen-pane view / active selection at K5 (previously
#    scrolled to row 13 / selecting L21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: C5:H5 already carry style s="2"; just populate their values ---
$ws.Range("C5:H5").Value = 5

# --- Row 5: new cell J5 (style s="6", matches J15/J26/I13/I26) ---
$ws.Range("J5").Value = 5
$ws.Range("J15").Copy()
$ws.Range("J5").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Row 9: new cell J9 + now-populated K9 (style s="4", matches I9/J21/K21) ---
$ws.Range("J9").Value = 5
$ws.Range("K9").Value = 5
$ws.Range("I9").Copy()
$ws.Range("J9:K9").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- View state: frozen pane now shows K5 as the active selection ---
$ws.Activate()
$ws.Range("K5").Select()
